# Generate Report for Handoff
# Updates the "Latest Handoff Date"/"Latest Handoff Datetime" timestamps
# on the Overview, zh-cn and de-de sheets to reflect a freshly regenerated
# handoff report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: column D holds "Latest Handoff Date" text values ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in @(7, 9, 10, 11, 13)) {
    $wsOverview.Range("D$r").Value = "2016-26-12 04:26:06"
}
foreach ($r in @(12, 14)) {
    $wsOverview.Range("D$r").Value = "2016-26-12 04:26:03"
}

# --- zh-cn sheet: column E holds "Latest Handoff Datetime" text values ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in @(7, 9, 10, 11, 12, 13, 14)) {
    $wsZhCn.Range("E$r").Value = "2016-03-12 04:26:03"
}

# --- de-de sheet: column E holds "Latest Handoff Datetime" text values ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in @(7, 9, 10, 11, 13)) {
    $wsDeDe.Range("E$r").Value = "2016-03-12 04:26:06"
}
